$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8 (existing rows 8-10 shift down to 9-11)
$ws.Rows.Item(8).Insert()

# Populate the new row with the "Energia.png" source entry
$ws.Range("A8").Value = "Energia.png"
$ws.Range("B8").Value = "https://opengameart.org/content/energy-icon"
$ws.Range("C8").Value = "CC-BY 4.0"

# The row insert does not move the hyperlink anchors automatically, so the
# hyperlinks that used to sit on B9/B10 are now stale (still anchored at
# B9/B10, but those rows now hold different data). Remove just those two.
foreach ($hl in $ws.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq "`$B`$9" -or $addr -eq "`$B`$10") {
        $hl.Delete()
    }
}

# Re-create the hyperlinks on their new (shifted) rows, plus the new one.
$ws.Hyperlinks.Add($ws.Range("B10"), "https://elthen.itch.io/2d-pixel-art-vegetable-monsters-sprite-pack")
$ws.Hyperlinks.Add($ws.Range("B11"), "https://free-game-assets.itch.io/night-city-street-2d-background-tiles")
$ws.Hyperlinks.Add($ws.Range("B8"), "https://opengameart.org/content/energy-icon")

# Match the existing hyperlink-cell formatting used elsewhere in column B.
$ws.Range("B8").Style = "Hyperlink"

# Update the selection to reflect where editing left off.
$ws.Range("C8").Select()
